# Update NATMI ligand-receptor edge statistics (Angpt2-Tie1) with recomputed TPM-derived values.
# Only columns G:T across data rows 2-21 change; columns A:F (cluster/gene labels, counts) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 13.40979833333333
$ws.Range("H2").Value = 40.229395
$ws.Range("I2").Value = 0.6868409202994065
$ws.Range("J2").Value = 0.6868409202994064
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 93.78975666666668
$ws.Range("N2").Value = 281.36927
$ws.Range("O2").Value = 0.997863063099077
$ws.Range("P2").Value = 0.9978630630990771
$ws.Range("Q2").Value = 1257.701722632406
$ws.Range("R2").Value = 11319.31550369165
$ws.Range("S2").Value = 0.6853731845917548
$ws.Range("T2").Value = 0.6853731845917548

# Row 3
$ws.Range("G3").Value = 13.40979833333333
$ws.Range("H3").Value = 40.229395
$ws.Range("I3").Value = 0.6868409202994065
$ws.Range("J3").Value = 0.6868409202994064
$ws.Range("K3").Value = 1.0
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.09680433333333333
$ws.Range("N3").Value = 0.290413
$ws.Range("O3").Value = 0.001029936231997873
$ws.Range("P3").Value = 0.001029936231997873
$ws.Range("Q3").Value = 1.298126587792778
$ws.Range("R3").Value = 11.683139290135
$ws.Range("S3").Value = 0.0007074023494351222
$ws.Range("T3").Value = 0.0007074023494351222

# Row 4
$ws.Range("G4").Value = 13.40979833333333
$ws.Range("H4").Value = 40.229395
$ws.Range("I4").Value = 0.6868409202994065
$ws.Range("J4").Value = 0.6868409202994064
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 0.05343666666666667
$ws.Range("N4").Value = 0.16031
$ws.Range("O4").Value = 0.0005685319780849309
$ws.Range("P4").Value = 0.000568531978084931
$ws.Range("Q4").Value = 0.7165749236055555
$ws.Range("R4").Value = 6.449174312449999
$ws.Range("S4").Value = 0.0003904910270474959
$ws.Range("T4").Value = 0.0003904910270474959

# Row 5
$ws.Range("G5").Value = 13.40979833333333
$ws.Range("H5").Value = 40.229395
$ws.Range("I5").Value = 0.6868409202994065
$ws.Range("J5").Value = 0.6868409202994064
$ws.Range("K5").Value = 2.0
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.050611
$ws.Range("N5").Value = 0.151833
$ws.Range("O5").Value = 0.0005384686908400556
$ws.Range("P5").Value = 0.0005384686908400557
$ws.Range("Q5").Value = 0.6786833034483333
$ws.Range("R5").Value = 6.108149731035
$ws.Range("S5").Value = 0.0003698423311690004
$ws.Range("T5").Value = 0.0003698423311690004

# Row 6
$ws.Range("G6").Value = 1.142193666666667
$ws.Range("H6").Value = 3.426581
$ws.Range("I6").Value = 0.05850239725256769
$ws.Range("J6").Value = 0.05850239725256768
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 93.78975666666668
$ws.Range("N6").Value = 281.36927
$ws.Range("O6").Value = 0.997863063099077
$ws.Range("P6").Value = 0.9978630630990771
$ws.Range("Q6").Value = 107.1260660628744
$ws.Range("R6").Value = 964.13459456587
$ws.Range("S6").Value = 0.05837738132108623
$ws.Range("T6").Value = 0.05837738132108622

# Row 7
$ws.Range("G7").Value = 1.142193666666667
$ws.Range("H7").Value = 3.426581
$ws.Range("I7").Value = 0.05850239725256769
$ws.Range("J7").Value = 0.05850239725256768
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.09680433333333333
$ws.Range("N7").Value = 0.290413
$ws.Range("O7").Value = 0.001029936231997873
$ws.Range("P7").Value = 0.001029936231997873
$ws.Range("Q7").Value = 0.1105692964392222
$ws.Range("R7").Value = 0.9951236679529998
$ws.Range("S7").Value = 0.0000602537385891523
$ws.Range("T7").Value = 0.0000602537385891523

# Row 8
$ws.Range("G8").Value = 1.142193666666667
$ws.Range("H8").Value = 3.426581
$ws.Range("I8").Value = 0.05850239725256769
$ws.Range("J8").Value = 0.05850239725256768
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 0.05343666666666667
$ws.Range("N8").Value = 0.16031
$ws.Range("O8").Value = 0.0005685319780849309
$ws.Range("P8").Value = 0.000568531978084931
$ws.Range("Q8").Value = 0.06103502223444444
$ws.Range("R8").Value = 0.54931520011
$ws.Range("S8").Value = 0.00003326048363271274
$ws.Range("T8").Value = 0.00003326048363271274

# Row 9
$ws.Range("G9").Value = 1.142193666666667
$ws.Range("H9").Value = 3.426581
$ws.Range("I9").Value = 0.05850239725256769
$ws.Range("J9").Value = 0.05850239725256768
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.050611
$ws.Range("N9").Value = 0.151833
$ws.Range("O9").Value = 0.0005384686908400556
$ws.Range("P9").Value = 0.0005384686908400557
$ws.Range("Q9").Value = 0.05780756366366666
$ws.Range("R9").Value = 0.5202680729729999
$ws.Range("S9").Value = 0.00003150170925959499
$ws.Range("T9").Value = 0.00003150170925959499

# Row 10
$ws.Range("G10").Value = 2.395418333333333
$ws.Range("H10").Value = 7.186254999999999
$ws.Range("I10").Value = 0.1226917282177923
$ws.Range("J10").Value = 0.1226917282177922
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 93.78975666666668
$ws.Range("N10").Value = 281.36927
$ws.Range("O10").Value = 0.997863063099077
$ws.Range("P10").Value = 0.9978630630990771
$ws.Range("Q10").Value = 224.6657025982056
$ws.Range("R10").Value = 2021.99132338385
$ws.Range("S10").Value = 0.1224295437363257
$ws.Range("T10").Value = 0.1224295437363256

# Row 11
$ws.Range("G11").Value = 2.395418333333333
$ws.Range("H11").Value = 7.186254999999999
$ws.Range("I11").Value = 0.1226917282177923
$ws.Range("J11").Value = 0.1226917282177922
$ws.Range("K11").Value = 1.0
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.09680433333333333
$ws.Range("N11").Value = 0.290413
$ws.Range("O11").Value = 0.001029936231997873
$ws.Range("P11").Value = 0.001029936231997873
$ws.Range("Q11").Value = 0.2318868748127778
$ws.Range("R11").Value = 2.086981873315
$ws.Range("S11").Value = 0.0001263646562579401
$ws.Range("T11").Value = 0.0001263646562579401

# Row 12
$ws.Range("G12").Value = 2.395418333333333
$ws.Range("H12").Value = 7.186254999999999
$ws.Range("I12").Value = 0.1226917282177923
$ws.Range("J12").Value = 0.1226917282177922
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 0.05343666666666667
$ws.Range("N12").Value = 0.16031
$ws.Range("O12").Value = 0.0005685319780849309
$ws.Range("P12").Value = 0.000568531978084931
$ws.Range("Q12").Value = 0.1280031710055556
$ws.Range("R12").Value = 1.15202853905
$ws.Range("S12").Value = 0.00006975417093832018
$ws.Range("T12").Value = 0.00006975417093832018

# Row 13
$ws.Range("G13").Value = 2.395418333333333
$ws.Range("H13").Value = 7.186254999999999
$ws.Range("I13").Value = 0.1226917282177923
$ws.Range("J13").Value = 0.1226917282177922
$ws.Range("K13").Value = 2.0
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.050611
$ws.Range("N13").Value = 0.151833
$ws.Range("O13").Value = 0.0005384686908400556
$ws.Range("P13").Value = 0.0005384686908400557
$ws.Range("Q13").Value = 0.1212345172683333
$ws.Range("R13").Value = 1.091110655415
$ws.Range("S13").Value = 0.00006606565427033851
$ws.Range("T13").Value = 0.00006606565427033851

# Row 14
$ws.Range("G14").Value = 1.306961666666667
$ws.Range("H14").Value = 3.920885
$ws.Range("I14").Value = 0.06694170423860808
$ws.Range("J14").Value = 0.06694170423860807
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 93.78975666666668
$ws.Range("N14").Value = 281.36927
$ws.Range("O14").Value = 0.997863063099077
$ws.Range("P14").Value = 0.9978630630990771
$ws.Range("Q14").Value = 122.5796166893278
$ws.Range("R14").Value = 1103.21655020395
$ws.Range("S14").Value = 0.06679865404060993
$ws.Range("T14").Value = 0.06679865404060992

# Row 15
$ws.Range("G15").Value = 1.306961666666667
$ws.Range("H15").Value = 3.920885
$ws.Range("I15").Value = 0.06694170423860808
$ws.Range("J15").Value = 0.06694170423860807
$ws.Range("K15").Value = 1.0
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.09680433333333333
$ws.Range("N15").Value = 0.290413
$ws.Range("O15").Value = 0.001029936231997873
$ws.Range("P15").Value = 0.001029936231997873
$ws.Range("Q15").Value = 0.1265195528338889
$ws.Range("R15").Value = 1.138675975505
$ws.Range("S15").Value = 0.00006894568662702806
$ws.Range("T15").Value = 0.00006894568662702806

# Row 16
$ws.Range("G16").Value = 1.306961666666667
$ws.Range("H16").Value = 3.920885
$ws.Range("I16").Value = 0.06694170423860808
$ws.Range("J16").Value = 0.06694170423860807
$ws.Range("K16").Value = 3.0
$ws.Range("L16").Value = 1.0
$ws.Range("M16").Value = 0.05343666666666667
$ws.Range("N16").Value = 0.16031
$ws.Range("O16").Value = 0.0005685319780849309
$ws.Range("P16").Value = 0.000568531978084931
$ws.Range("Q16").Value = 0.06983967492777779
$ws.Range("R16").Value = 0.6285570743500001
$ws.Range("S16").Value = 0.00003805849952715225
$ws.Range("T16").Value = 0.00003805849952715225

# Row 17
$ws.Range("G17").Value = 1.306961666666667
$ws.Range("H17").Value = 3.920885
$ws.Range("I17").Value = 0.06694170423860808
$ws.Range("J17").Value = 0.06694170423860807
$ws.Range("K17").Value = 2.0
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.050611
$ws.Range("N17").Value = 0.151833
$ws.Range("O17").Value = 0.0005384686908400556
$ws.Range("P17").Value = 0.0005384686908400557
$ws.Range("Q17").Value = 0.06614663691166667
$ws.Range("R17").Value = 0.595319732205
$ws.Range("S17").Value = 0.00003604601184396549
$ws.Range("T17").Value = 0.00003604601184396549

# Row 18
$ws.Range("G18").Value = 1.269506
$ws.Range("H18").Value = 3.808518
$ws.Range("I18").Value = 0.06502324999162565
$ws.Range("J18").Value = 0.06502324999162565
$ws.Range("K18").Value = 3.0
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 93.78975666666668
$ws.Range("N18").Value = 281.36927
$ws.Range("O18").Value = 0.997863063099077
$ws.Range("P18").Value = 0.9978630630990771
$ws.Range("Q18").Value = 119.0666588268734
$ws.Range("R18").Value = 1071.59992944186
$ws.Range("S18").Value = 0.0648842994093006
$ws.Range("T18").Value = 0.06488429940930061

# Row 19
$ws.Range("G19").Value = 1.269506
$ws.Range("H19").Value = 3.808518
$ws.Range("I19").Value = 0.06502324999162565
$ws.Range("J19").Value = 0.06502324999162565
$ws.Range("K19").Value = 1.0
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.09680433333333333
$ws.Range("N19").Value = 0.290413
$ws.Range("O19").Value = 0.001029936231997873
$ws.Range("P19").Value = 0.001029936231997873
$ws.Range("Q19").Value = 0.1228936819926667
$ws.Range("R19").Value = 1.106043137934
$ws.Range("S19").Value = 0.00006696980108863065
$ws.Range("T19").Value = 0.00006696980108863067

# Row 20
$ws.Range("G20").Value = 1.269506
$ws.Range("H20").Value = 3.808518
$ws.Range("I20").Value = 0.06502324999162565
$ws.Range("J20").Value = 0.06502324999162565
$ws.Range("K20").Value = 3.0
$ws.Range("L20").Value = 1.0
$ws.Range("M20").Value = 0.05343666666666667
$ws.Range("N20").Value = 0.16031
$ws.Range("O20").Value = 0.0005685319780849309
$ws.Range("P20").Value = 0.000568531978084931
$ws.Range("Q20").Value = 0.06783816895333333
$ws.Range("R20").Value = 0.61054352058
$ws.Range("S20").Value = 0.0000369677969392499
$ws.Range("T20").Value = 0.0000369677969392499

# Row 21
$ws.Range("G21").Value = 1.269506
$ws.Range("H21").Value = 3.808518
$ws.Range("I21").Value = 0.06502324999162565
$ws.Range("J21").Value = 0.06502324999162565
$ws.Range("K21").Value = 2.0
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.050611
$ws.Range("N21").Value = 0.151833
$ws.Range("O21").Value = 0.0005384686908400556
$ws.Range("P21").Value = 0.0005384686908400557
$ws.Range("Q21").Value = 0.064250968166
$ws.Range("R21").Value = 0.578258713049999
$ws.Range("S21").Value = 0.00003501298429715632
$ws.Range("T21").Value = 0.00003501298429715632
